# Update the cached "datetimeFigureOut" date field text from 5/20/2025 to
# 5/30/2025 across the slide master, every slide layout, and the notes
# master (mirrors PowerPoint silently refreshing the auto-date placeholder
# text it had cached the last time the deck was saved).

$p = $ppt.ActivePresentation

$oldDate = "5/20/2025"
$newDate = "5/30/2025"

function Update-DatePlaceholder($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $sh = $container.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            $isDatePlaceholder = $false
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
            if ($isDatePlaceholder) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master's own Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master

# Every slide layout hanging off the master.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout
}

# Notes master's Date Placeholder.
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster
